$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3336.2666
$ws.Range("J17").Value = 3336.2666
$ws.Range("L17").Value = 10008.7998
$ws.Range("N17").Value = -10344.7998
$ws.Range("H18").Value = 1174.8
$ws.Range("I18").Value = 1174.8
$ws.Range("K18").Value = 1174.8
$ws.Range("M18").Value = -890.8
$ws.Range("H32").Value = 4690.6665
$ws.Range("I32").Value = 3750
$ws.Range("J32").Value = 4878.8
$ws.Range("K32").Value = 3750
$ws.Range("L32").Value = 4878.8
$ws.Range("M32").Value = -3424
$ws.Range("N32").Value = -5530.8
$ws.Range("H38").Value = 1419.6
$ws.Range("J38").Value = 1799
$ws.Range("L38").Value = 5397
$ws.Range("N38").Value = -6141
$ws.Range("H137").Value = 6336.6665
$ws.Range("I137").Value = 1754.625
$ws.Range("J137").Value = 10002.3
$ws.Range("K137").Value = 5263.875
$ws.Range("L137").Value = 30006.9
$ws.Range("M137").Value = -2713.875
$ws.Range("N137").Value = -35106.89999999999
$ws.Range("H138").Value = 3968.9285
$ws.Range("I138").Value = 2523.25
$ws.Range("J138").Value = 5896.5
$ws.Range("K138").Value = 7569.75
$ws.Range("L138").Value = 17689.5
$ws.Range("M138").Value = -2429.75
$ws.Range("N138").Value = -27969.5

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1362.5927
$ws.Range("I2").Value = 1280.8572
$ws.Range("K2").Value = 1280.8572
$ws.Range("M2").Value = -1167.8572
$ws.Range("H61").Value = 2957.2307
$ws.Range("I61").Value = 1953.6666
$ws.Range("K61").Value = 1953.6666
$ws.Range("M61").Value = -1741.6666
$ws.Range("H62").Value = 53637.5
$ws.Range("J62").Value = 53637.5
$ws.Range("L62").Value = 53637.5
$ws.Range("N62").Value = -54885.5
$ws.Range("H65").Value = 53637.5
$ws.Range("J65").Value = 53637.5
$ws.Range("L65").Value = 160912.5
$ws.Range("N65").Value = -167152.5
$ws.Range("H74").Value = 3773.9092
$ws.Range("I74").Value = 1127.0454
$ws.Range("J74").Value = 9067.637000000001
$ws.Range("K74").Value = 1127.0454
$ws.Range("L74").Value = 9067.637000000001
$ws.Range("M74").Value = -253.0454
$ws.Range("N74").Value = -10815.637
$ws.Range("H77").Value = 3773.9092
$ws.Range("I77").Value = 1127.0454
$ws.Range("J77").Value = 9067.637000000001
$ws.Range("K77").Value = 5635.227
$ws.Range("L77").Value = 45338.185
$ws.Range("M77").Value = -1267.227
$ws.Range("N77").Value = -54074.185
$ws.Range("H116").Value = 1362.5927
$ws.Range("I116").Value = 1280.8572
$ws.Range("K116").Value = 1280.8572
$ws.Range("M116").Value = 1013.1428
$ws.Range("H122").Value = 2632.7576
$ws.Range("I122").Value = 1982.9615
$ws.Range("K122").Value = 5948.8845
$ws.Range("M122").Value = -3498.8845
$ws.Range("H132").Value = 2725.25
$ws.Range("I132").Value = 1574.3043
$ws.Range("K132").Value = 4722.9129
$ws.Range("M132").Value = -2192.9129
$ws.Range("H136").Value = 2957.2307
$ws.Range("I136").Value = 1953.6666
$ws.Range("K136").Value = 5860.9998
$ws.Range("M136").Value = -3310.9998

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1362.5927
$ws.Range("I3").Value = 1280.8572
$ws.Range("K3").Value = 1280.8572
$ws.Range("M3").Value = -1166.8572
$ws.Range("H86").Value = 8686.429
$ws.Range("I86").Value = 2599.75
$ws.Range("K86").Value = 2599.75
$ws.Range("M86").Value = -1476.75
$ws.Range("H89").Value = 8686.429
$ws.Range("I89").Value = 2599.75
$ws.Range("K89").Value = 12998.75
$ws.Range("M89").Value = -7382.75
$ws.Range("H105").Value = 5250.5483
$ws.Range("I105").Value = 3911.9092
$ws.Range("K105").Value = 3911.9092
$ws.Range("M105").Value = -2164.9092
$ws.Range("H107").Value = 1793.6538
$ws.Range("I107").Value = 1872.9524
$ws.Range("J107").Value = 1460.6
$ws.Range("K107").Value = 1872.9524
$ws.Range("L107").Value = 1460.6
$ws.Range("M107").Value = 47.0476000000001
$ws.Range("N107").Value = -5300.6
$ws.Range("H134").Value = 3312.6667
$ws.Range("I134").Value = 1497.2222
$ws.Range("J134").Value = 8759
$ws.Range("K134").Value = 4491.6666
$ws.Range("L134").Value = 26277
$ws.Range("M134").Value = -1956.6666
$ws.Range("N134").Value = -31347

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6634.1665
$ws.Range("I31").Value = 3232.8635
$ws.Range("J31").Value = 11979.071
$ws.Range("K31").Value = 3232.8635
$ws.Range("L31").Value = 11979.071
$ws.Range("M31").Value = -2937.8635
$ws.Range("N31").Value = -12569.071
$ws.Range("H34").Value = 6634.1665
$ws.Range("I34").Value = 3232.8635
$ws.Range("J34").Value = 11979.071
$ws.Range("K34").Value = 3232.8635
$ws.Range("L34").Value = 11979.071
$ws.Range("M34").Value = -3030.8635
$ws.Range("N34").Value = -12383.071
$ws.Range("H97").Value = 149608.5
$ws.Range("I97").Value = 60000
$ws.Range("K97").Value = 60000
$ws.Range("M97").Value = -59009
$ws.Range("H132").Value = 2958.5405
$ws.Range("I132").Value = 2460.7942
$ws.Range("K132").Value = 7382.382599999999
$ws.Range("M132").Value = -4852.382599999999

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 186.7
$ws.Range("I33").Value = 189.75
$ws.Range("J33").Value = 174.5
$ws.Range("K33").Value = 1138.5
$ws.Range("L33").Value = 1047
$ws.Range("M33").Value = -855.5
$ws.Range("N33").Value = -1613
$ws.Range("H122").Value = 661.9048
$ws.Range("J122").Value = 619.44446
$ws.Range("L122").Value = 5575.00014
$ws.Range("N122").Value = -10475.00014
$ws.Range("H131").Value = 1300631.5
$ws.Range("I131").Value = 1149
$ws.Range("K131").Value = 3447
$ws.Range("M131").Value = 1593

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4366.25
$ws.Range("I80").Value = 3504.6667
$ws.Range("J80").Value = 5071.1816
$ws.Range("K80").Value = 3504.6667
$ws.Range("L80").Value = 5071.1816
$ws.Range("M80").Value = -2506.6667
$ws.Range("N80").Value = -7067.1816
$ws.Range("H83").Value = 4366.25
$ws.Range("I83").Value = 3504.6667
$ws.Range("J83").Value = 5071.1816
$ws.Range("K83").Value = 17523.3335
$ws.Range("L83").Value = 25355.908
$ws.Range("M83").Value = -12531.3335
$ws.Range("N83").Value = -35339.908
$ws.Range("H132").Value = 2849.027
$ws.Range("I132").Value = 2278.0715
$ws.Range("J132").Value = 4625.3335
$ws.Range("K132").Value = 6834.2145
$ws.Range("L132").Value = 13876.0005
$ws.Range("M132").Value = -4304.2145
$ws.Range("N132").Value = -18936.0005

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 11499.934
$ws.Range("J20").Value = 21750
$ws.Range("L20").Value = 21750
$ws.Range("N20").Value = -22202
$ws.Range("H22").Value = 1132.3871
$ws.Range("I22").Value = 1068.5883
$ws.Range("J22").Value = 1209.8572
$ws.Range("K22").Value = 1068.5883
$ws.Range("L22").Value = 1209.8572
$ws.Range("M22").Value = -773.5882999999999
$ws.Range("N22").Value = -1799.8572
$ws.Range("H27").Value = 1132.3871
$ws.Range("I27").Value = 1068.5883
$ws.Range("J27").Value = 1209.8572
$ws.Range("K27").Value = 1068.5883
$ws.Range("L27").Value = 1209.8572
$ws.Range("M27").Value = -961.5882999999999
$ws.Range("N27").Value = -1423.8572
$ws.Range("H36").Value = 65554
$ws.Range("J36").Value = 65554
$ws.Range("L36").Value = 65554
$ws.Range("N36").Value = -66678
$ws.Range("H40").Value = 10777
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("H46").Value = 3551.4375
$ws.Range("I46").Value = 1734.375
$ws.Range("J46").Value = 5368.5
$ws.Range("K46").Value = 1734.375
$ws.Range("L46").Value = 5368.5
$ws.Range("M46").Value = -1546.375
$ws.Range("N46").Value = -5744.5
$ws.Range("H55").Value = 1315.9333
$ws.Range("I55").Value = 677.1111
$ws.Range("J55").Value = 2274.1667
$ws.Range("K55").Value = 677.1111
$ws.Range("L55").Value = 2274.1667
$ws.Range("M55").Value = -504.1111
$ws.Range("N55").Value = -2620.1667
$ws.Range("H132").Value = 6757.659
$ws.Range("I132").Value = 3498.682
$ws.Range("J132").Value = 10016.637
$ws.Range("K132").Value = 10496.046
$ws.Range("L132").Value = 30049.911
$ws.Range("M132").Value = -7966.045999999998
$ws.Range("N132").Value = -35109.911
$ws.Range("H133").Value = 95000
$ws.Range("J133").Value = 95000
$ws.Range("L133").Value = 95000
$ws.Range("N133").Value = -100060
$ws.Range("H136").Value = 7024.2983
$ws.Range("I136").Value = 4065.9524
$ws.Range("K136").Value = 12197.8572
$ws.Range("M136").Value = -9647.8572

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H13").Value = 2300
$ws.Range("I13").Value = 2300
$ws.Range("K13").Value = 2300
$ws.Range("M13").Value = -2160
$ws.Range("H81").Value = 1807.2727
$ws.Range("J81").Value = 3933
$ws.Range("L81").Value = 7866
$ws.Range("N81").Value = -9988
$ws.Range("H84").Value = 1807.2727
$ws.Range("J84").Value = 3933
$ws.Range("L84").Value = 39330
$ws.Range("N84").Value = -49938
$ws.Range("H107").Value = 1426.5769
$ws.Range("I107").Value = 1284.5238
$ws.Range("K107").Value = 3853.5714
$ws.Range("M107").Value = -1933.5714
$ws.Range("H110").Value = 163300
$ws.Range("J110").Value = 163300
$ws.Range("L110").Value = 163300
$ws.Range("N110").Value = -171480
$ws.Range("H132").Value = 3920.239
$ws.Range("I132").Value = 2859.9697
$ws.Range("J132").Value = 6611.6924
$ws.Range("K132").Value = 8579.909100000001
$ws.Range("L132").Value = 19835.0772
$ws.Range("M132").Value = -6049.909100000001
$ws.Range("N132").Value = -24895.0772
